$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-breaking space used after "Moeen Ali" in the existing rows of this sheet.
$nbsp = [char]0x00A0
$batsman = "Moeen Ali" + $nbsp

# The sheet currently holds a header row + 3 match rows (rows 1-4).
# Append the same 3 match rows again (rows 5-7), duplicating rows 2-4,
# exactly as the source data feed re-emitted them.
$rowsData = @(
    @(" Dubai (DSC)", " October 25 2020", "Super Kings won by 8 wickets (with 8 balls remaining)", "Royal Challengers Bangalore", "Chennai Super Kings", $batsman, "1", "2", "0", "0", "50.00"),
    @(" Abu Dhabi", " November 06 2020", "Sunrisers won by 6 wickets (with 2 balls remaining)", "Royal Challengers Bangalore", "Sunrisers Hyderabad", $batsman, "0", "1", "0", "0", "0.00"),
    @(" Dubai (DSC)", " October 05 2020", "Capitals won by 59 runs", "Royal Challengers Bangalore", "Delhi Capitals", $batsman, "11", "13", "1", "0", "84.61")
)

$startRow = 5
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $rowValues = $rowsData[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $value = $rowValues[$c - 1]
        $cell = $ws.Cells.Item($r, $c)
        if ($c -ge 7) {
            # Columns G:K hold numeric-looking values that must stay text,
            # matching the existing rows (t="str") - use a leading
            # apostrophe so Excel keeps them as text instead of numbers,
            # then restore the default "Normal" style so no quote-prefix
            # formatting artifact is left behind on the cell.
            $cell.Value = "'" + $value
            $cell.Style = "Normal"
        } else {
            $cell.Value = $value
        }
    }
}
